# Project_Demo_Form.docx revision ("All Done! Please use this commit for
# grading!"):
#
# Three wording tweaks -- each member's contribution blurb gets the word
# "initial " inserted to clarify that the *first* ER Diagram / SQL
# conversion is being referred to:
#
#   "Haojia Yu: Suggestions on the ER Diagram creation."
#     -> "Haojia Yu: Suggestions on the ER Diagram initial creation."
#
#   "Yuxiang Huang: Convert the ER Diagram to sql commands."
#     -> "Yuxiang Huang: Convert the initial ER Diagram to sql commands."
#
#   "Haojia Yu: Convert the ER Diagram to sql commands."
#     -> "Haojia Yu: Convert the initial ER Diagram to sql commands."
#
# (Everything else in the upstream revision -- the reordering of runs, the
# removal of the spell-check <w:proofErr> squiggly-line markers, the
# collapsing of several runs into one -- is purely a cosmetic side effect
# of Word re-serialising the paragraphs it touched; none of it changes the
# document's visible text, formatting, or paragraph structure, so the
# three Find/Replace edits below are the whole of the substantive change.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Haojia Yu: Suggestions on the ER Diagram creation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Haojia Yu: Suggestions on the ER Diagram initial creation.", 2) | Out-Null

$d.Content.Find.Execute(
    "Yuxiang Huang: Convert the ER Diagram to sql commands.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Yuxiang Huang: Convert the initial ER Diagram to sql commands.", 2) | Out-Null

$d.Content.Find.Execute(
    "Haojia Yu: Convert the ER Diagram to sql commands.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Haojia Yu: Convert the initial ER Diagram to sql commands.", 2) | Out-Null
